$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous table content entirely before laying out the new grant table.
$ws.Range("A1:E13").ClearContents()

# Row 1: header
$ws.Range("A1").Value = 'order'
$ws.Range("B1").Value = 'what'
$ws.Range("C1").Value = 'when'
$ws.Range("D1").Value = 'with'
$ws.Range("E1").Value = 'where'
$ws.Range("F1").Value = 'why'

# Grant #1 (rows 2-5)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'National Institute of Health: National Institute on Aging'
$ws.Range("C2").Value = 'Sept 2000 - Sept 2005'
$ws.Range("D2").Value = 'R01 Grant: Cache County Family-based Cohort Study on Aging.'
$ws.Range("E2").Value = 'Utah State University'
$ws.Range("F2").Value = 'Pis: Drs John Breitner, Kathy Welch-Bohmer'
$ws.Range("F3").Value = 'Roll: data management and analysis'
$ws.Range("F4").Value = 'Amount \$1,999,400(original)'
$ws.Range("F5").Value = 'Extended multiple times'

# Grant #2 (rows 6-10)
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 'National Institute of Health: National Institute on Aging'
$ws.Range("C6").Value = 'Sept. 2002 – Sept. 2013'
$ws.Range("D6").Value = ' R01 Grant: Progression of Dementia, A Population Study. '
$ws.Range("E6").Value = 'Utah State University'
$ws.Range("F6").Value = 'PIs: Dr. Joann Tschanz and  Dr. Constantine G. Lyketsos  '
$ws.Range("F7").Value = 'Roll: data management and analysis'
$ws.Range("F8").Value = 'DUNS ID: 072983455 (original)'
$ws.Range("F9").Value = 'Amount: \$2,787,792 (original)'
$ws.Range("F10").Value = 'Extended multiple times'

# Grant #3 (rows 11-14)
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = 'National Institute of Health: National Institute on Aging'
$ws.Range("C11").Value = 'July 2008 - June 2011'
$ws.Range("D11").Value = 'R01 Grant: Lifespan Stressors and Alzheimer’s Disease: The Cache County Study.'
$ws.Range("E11").Value = 'Utah State University'
$ws.Range("F11").Value = 'Pis: Dr. Maria Norton'
$ws.Range("F12").Value = 'Co-PI: Dr. Joann Tschanz'
$ws.Range("F13").Value = 'Roll: data management and analysis'
$ws.Range("F14").Value = 'Amount:  \$970,549 (original)'

# Grant #4 (rows 19-23)
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = 'National Science Foundation: Division Of Research On Learning'
$ws.Range("C19").Value = 'Sept. 2019 - Feb. 2022'
$ws.Range("D19").Value = 'Research on the Development of An Assessment to Measure Kindergarten Children''s Abilities to Reason Computationally With Mathematical Problem-Solving Skills'
$ws.Range("E19").Value = 'Utah State University'
$ws.Range("F19").Value = 'Program: STEM + Computing (STEM+C) Part'
$ws.Range("F20").Value = 'PI: Jody Clarke Midura'
$ws.Range("F21").Value = 'Co-PI: Victor Raymond Lee, Jessica Shumway'
$ws.Range("F22").Value = 'DUNS ID: 072983455'
$ws.Range("F23").Value = 'Amount: \$1,120,807'

# Currency number format on the last amount cell (F23)
$ws.Range("F23").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# Column widths
$ws.Columns.Item(2).ColumnWidth = 37.86328125
$ws.Columns.Item(3).ColumnWidth = 20.9296875
$ws.Columns.Item(4).ColumnWidth = 33.73046875
$ws.Columns.Item(5).ColumnWidth = 27.06640625
$ws.Columns.Item(6).ColumnWidth = 39.59765625

# Final selection, matching the saved cursor position
$ws.Range("F15").Select()
